$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 120 (pushes old rows 120:195 down to 121:196,
# carrying their formatting/styles with them).
$ws.Rows("120:120").Insert()

# Populate the newly inserted row 120 with a new weekly price record
# (template values copied from the surrounding Cebollín entries, with a
# new date and volume/price figures).
$ws.Cells.Item(120, 1).Value = 4
$ws.Cells.Item(120, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(120, 3).Value = "Los Lagos"
$ws.Cells.Item(120, 4).Value = 44529
$ws.Cells.Item(120, 5).Value = 10
$ws.Cells.Item(120, 6).Value = 100112037
$ws.Cells.Item(120, 7).Value = "Cebollín"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 90
$ws.Cells.Item(120, 11).Value = 6000
$ws.Cells.Item(120, 12).Value = 6000
$ws.Cells.Item(120, 13).Value = 6000
$ws.Cells.Item(120, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(120, 15).Value = "Región Metropolitana"
$ws.Cells.Item(120, 16).Value = 167
$ws.Cells.Item(120, 17).Value = 36
$ws.Cells.Item(120, 18).Value = "Hortaliza"
